# Auto-generated edit script applying numeric corrections per the commit diff
# (values recomputed by the scheduled Gilgamesh profits runner).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 255
$ws.Range("I2").Value = 217.3077
$ws.Range("K2").Value = 217.3077
$ws.Range("M2").Value = -104.3077
$ws.Range("H9").Value = 710815.4
$ws.Range("I9").Value = 1066148.6
$ws.Range("K9").Value = 1066148.6
$ws.Range("M9").Value = -1065979.6
$ws.Range("H86").Value = 6563
$ws.Range("I86").Value = 5513.273
$ws.Range("K86").Value = 5513.273
$ws.Range("M86").Value = -4390.273
$ws.Range("H89").Value = 6563
$ws.Range("I89").Value = 5513.273
$ws.Range("K89").Value = 27566.365
$ws.Range("M89").Value = -21950.365
$ws.Range("H100").Value = 6033.5
$ws.Range("I100").Value = 5800.294
$ws.Range("J100").Value = 9998
$ws.Range("K100").Value = 5800.294
$ws.Range("L100").Value = 9998
$ws.Range("M100").Value = -5259.294
$ws.Range("N100").Value = -11080
$ws.Range("H107").Value = 570.7692
$ws.Range("I107").Value = 570.7692
$ws.Range("K107").Value = 570.7692
$ws.Range("M107").Value = 1349.2308
$ws.Range("H115").Value = 1844.6
$ws.Range("J115").Value = 1998
$ws.Range("L115").Value = 5994
$ws.Range("N115").Value = -9128
$ws.Range("H137").Value = 4208.4287
$ws.Range("I137").Value = 1501.3684
$ws.Range("J137").Value = 7423.0625
$ws.Range("K137").Value = 4504.1052
$ws.Range("L137").Value = 22269.1875
$ws.Range("M137").Value = -1954.1052
$ws.Range("N137").Value = -27369.1875
$ws.Range("H138").Value = 4170.615
$ws.Range("J138").Value = 7080
$ws.Range("L138").Value = 21240
$ws.Range("N138").Value = -31520

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4806.6665
$ws.Range("I2").Value = 2662.1538
$ws.Range("K2").Value = 2662.1538
$ws.Range("M2").Value = -2549.1538
$ws.Range("H45").Value = 170404
$ws.Range("I45").Value = 170404
$ws.Range("K45").Value = 170404
$ws.Range("M45").Value = -170027
$ws.Range("H116").Value = 4806.6665
$ws.Range("I116").Value = 2662.1538
$ws.Range("K116").Value = 2662.1538
$ws.Range("M116").Value = -368.1538

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4806.6665
$ws.Range("I3").Value = 2662.1538
$ws.Range("K3").Value = 2662.1538
$ws.Range("M3").Value = -2548.1538
$ws.Range("H20").Value = 19382902
$ws.Range("I20").Value = 23812982
$ws.Range("J20").Value = 1304.625
$ws.Range("K20").Value = 23812982
$ws.Range("L20").Value = 1304.625
$ws.Range("M20").Value = -23812735
$ws.Range("N20").Value = -1798.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 62903
$ws.Range("J9").Value = 62903
$ws.Range("L9").Value = 62903
$ws.Range("N9").Value = -63239
$ws.Range("H16").Value = 83334160
$ws.Range("I16").Value = 83334160
$ws.Range("K16").Value = 83334160
$ws.Range("M16").Value = -83333873
$ws.Range("H53").Value = 71666.664
$ws.Range("J53").Value = 71666.664
$ws.Range("L53").Value = 71666.664
$ws.Range("N53").Value = -72880.664
$ws.Range("H105").Value = 1710.381
$ws.Range("J105").Value = 5000
$ws.Range("L105").Value = 5000
$ws.Range("N105").Value = -8494
$ws.Range("H113").Value = 83334160
$ws.Range("I113").Value = 83334160
$ws.Range("K113").Value = 83334160
$ws.Range("M113").Value = -83331990
$ws.Range("H132").Value = 2361.0334
$ws.Range("I132").Value = 2028.76
$ws.Range("K132").Value = 6086.28
$ws.Range("M132").Value = -3556.28
$ws.Range("H133").Value = 80000
$ws.Range("J133").Value = 80000
$ws.Range("L133").Value = 80000
$ws.Range("N133").Value = -85060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 316.33334
$ws.Range("I7").Value = 316.33334
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 949.0000200000001
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -837.0000200000001
$ws.Range("N7").ClearContents()
$ws.Range("H56").Value = 7806
$ws.Range("I56").Value = 7806
$ws.Range("K56").Value = 7806
$ws.Range("M56").Value = -7276
$ws.Range("H107").Value = 4196.7188
$ws.Range("I107").Value = 691.2222
$ws.Range("J107").Value = 5568.4346
$ws.Range("K107").Value = 2073.6666
$ws.Range("L107").Value = 16705.3038
$ws.Range("M107").Value = -153.6666
$ws.Range("N107").Value = -20545.3038
$ws.Range("H121").Value = 3893094.8
$ws.Range("I121").Value = 14300609
$ws.Range("K121").Value = 42901827
$ws.Range("M121").Value = -42900517
$ws.Range("H122").Value = 2796.75
$ws.Range("J122").Value = 2796.75
$ws.Range("L122").Value = 25170.75
$ws.Range("N122").Value = -30070.75
$ws.Range("H131").Value = 2105016.5
$ws.Range("I131").Value = 10602.429
$ws.Range("J131").Value = 2803154.5
$ws.Range("K131").Value = 31807.287
$ws.Range("L131").Value = 8409463.5
$ws.Range("M131").Value = -26767.287
$ws.Range("N131").Value = -8419543.5
$ws.Range("H132").Value = 5930.8613
$ws.Range("J132").Value = 5894.96
$ws.Range("L132").Value = 53054.64
$ws.Range("N132").Value = -58114.64

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 43481180
$ws.Range("I80").Value = 66668188
$ws.Range("K80").Value = 66668188
$ws.Range("M80").Value = -66667190
$ws.Range("H83").Value = 43481180
$ws.Range("I83").Value = 66668188
$ws.Range("K83").Value = 333340940
$ws.Range("M83").Value = -333335948
$ws.Range("H122").Value = 37041190
$ws.Range("I122").Value = 58825360
$ws.Range("K122").Value = 176476080
$ws.Range("M122").Value = -176473630
$ws.Range("H139").Value = 100000
$ws.Range("J139").Value = 100000
$ws.Range("L139").Value = 100000
$ws.Range("N139").Value = -110280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1777.3334
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H132").Value = 8865.700000000001
$ws.Range("J132").Value = 7308.091
$ws.Range("L132").Value = 21924.273
$ws.Range("N132").Value = -26984.273
$ws.Range("H136").Value = 5657.1113
$ws.Range("I136").Value = 4251.3335
$ws.Range("J136").Value = 8468.666999999999
$ws.Range("K136").Value = 12754.0005
$ws.Range("L136").Value = 25406.001
$ws.Range("M136").Value = -10204.0005
$ws.Range("N136").Value = -30506.001
$ws.Range("H140").Value = 106110.18
$ws.Range("J140").Value = 106110.18
$ws.Range("L140").Value = 106110.18
$ws.Range("N140").Value = -116470.18

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5228.2856
$ws.Range("I81").Value = 5041.6924
$ws.Range("J81").Value = 5338.5454
$ws.Range("K81").Value = 10083.3848
$ws.Range("L81").Value = 10677.0908
$ws.Range("M81").Value = -9022.3848
$ws.Range("N81").Value = -12799.0908
$ws.Range("H84").Value = 5228.2856
$ws.Range("I84").Value = 5041.6924
$ws.Range("J84").Value = 5338.5454
$ws.Range("K84").Value = 50416.924
$ws.Range("L84").Value = 53385.454
$ws.Range("M84").Value = -45112.924
$ws.Range("N84").Value = -63993.454
$ws.Range("H107").Value = 2321
$ws.Range("I107").Value = 784.6
$ws.Range("K107").Value = 2353.8
$ws.Range("M107").Value = -433.8000000000002
$ws.Range("H132").Value = 3587.8333
$ws.Range("I132").Value = 3368.1316
$ws.Range("K132").Value = 10104.3948
$ws.Range("M132").Value = -7574.3948
$ws.Range("H136").Value = 6966.6
$ws.Range("I136").Value = 5718.25
$ws.Range("J136").Value = 11960
$ws.Range("K136").Value = 17154.75
$ws.Range("L136").Value = 35880
$ws.Range("M136").Value = -14604.75
$ws.Range("N136").Value = -40980
